$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Expense"
$ws.Range("C1").Value = "Amount"

# --- Data rows (expense log) ---
$ws.Range("A2").Value = 45569
$ws.Range("B2").Value = "Food"
$ws.Range("C2").Value = 10

$ws.Range("A3").Value = 45570
$ws.Range("B3").Value = "Food"
$ws.Range("C3").Value = 20

$ws.Range("A4").Value = 45571
$ws.Range("B4").Value = "Gas"
$ws.Range("C4").Value = 40

$ws.Range("A5").Value = 45572
$ws.Range("B5").Value = "Water Bill"
$ws.Range("C5").Value = 30

$ws.Range("A6").Value = 45573
$ws.Range("B6").Value = "Electric Bill"
$ws.Range("C6").Value = 80

# --- Date number formatting (built-in format 14) applied to A2, then
#     copied down so every date cell shares the same style record ---
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.3
$ws.Columns.Item(2).ColumnWidth = 14.8

# --- Selection: whole column C active, matches last UI state in file ---
$ws.Columns.Item(3).Select() | Out-Null
